$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Statut column: "RA" -> "NA"
$ws.Range("E2").Value = "NA"

# "Date dernière adhésion" column: clear the stored "2024" text value
$ws.Range("F2").Value = $null

# Montant adhésion année n: 0 -> 10
$ws.Range("M2").Value = 10

# Total année n: 0 -> 10
$ws.Range("P2").Value = 10
